$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- Row 21 ----
Set-TextValue $ws.Cells.Item(21,1) "2024-07-01"
$ws.Cells.Item(21,2).Value = "신한글로벌액티브리츠"
$ws.Cells.Item(21,3).Value = "코스피"
$ws.Cells.Item(21,4).Value = 700.0000199999999
$ws.Cells.Item(21,5).Value = "신한"
$ws.Cells.Item(21,6).Value = 350.00001
$ws.Cells.Item(21,7).Value = "-"
$ws.Cells.Item(21,8).Value = "-"
$ws.Cells.Item(21,9).Value = "-"
$ws.Cells.Item(21,10).Value = "-"
$ws.Cells.Item(21,11).Value = "공동대표"
$ws.Cells.Item(21,12).Value = "-"
$ws.Cells.Item(21,13).Value = 3000
$ws.Cells.Item(21,14).Value = 50
Set-TextValue $ws.Cells.Item(21,15) "2024-06-13"
Set-TextValue $ws.Cells.Item(21,16) "2024-06-18"
$ws.Cells.Item(21,17).Value = 32666668

# ---- Row 22 ----
Set-TextValue $ws.Cells.Item(22,1) "2024-07-01"
$ws.Cells.Item(22,2).Value = "신한글로벌액티브리츠"
$ws.Cells.Item(22,3).Value = "코스피"
$ws.Cells.Item(22,4).Value = 700.0000199999999
$ws.Cells.Item(22,5).Value = "한국"
$ws.Cells.Item(22,6).Value = 350.00001
$ws.Cells.Item(22,7).Value = "-"
$ws.Cells.Item(22,8).Value = "-"
$ws.Cells.Item(22,9).Value = "-"
$ws.Cells.Item(22,10).Value = "-"
$ws.Cells.Item(22,11).Value = "공동대표"
$ws.Cells.Item(22,12).Value = "-"
$ws.Cells.Item(22,13).Value = 3000
$ws.Cells.Item(22,14).Value = 50
Set-TextValue $ws.Cells.Item(22,15) "2024-06-13"
Set-TextValue $ws.Cells.Item(22,16) "2024-06-18"
$ws.Cells.Item(22,17).Value = 32666668
